$wb = $excel.ActiveWorkbook

# --- Logs sheet: append new row 18 ---
$ws = $wb.Worksheets.Item("Logs")
$row = 18
$ws.Cells.Item($row, 1).Value = "Kun jij dit even regelen?"
$ws.Cells.Item($row, 2).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item($row, 3).Value = "Testmail #1: Kun jij dit even regelen?"
$ws.Cells.Item($row, 4).Value = "Planning / Afspraak"
$ws.Cells.Item($row, 5).Value = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
$ws.Cells.Item($row, 6).Value = "2025-08-05 17:40:10"
$ws.Cells.Item($row, 7).Value = "Ja"
$ws.Cells.Item($row, 8).Value = "Ja"
$ws.Cells.Item($row, 9).Value = "Nee"
$ws.Cells.Item($row, 10).Value = "Nee"

# --- Extend conditional formatting ranges to include the new row ---
$cols = @("D", "G", "H", "I", "J")
foreach ($col in $cols) {
    $oldRange = $ws.Range($col + "2:" + $col + "17")
    $newRange = $ws.Range($col + "2:" + $col + "18")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Dashboard sheet: bump the "Planning / Afspraak" count from 11 to 12 ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 12
